$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly price record as row 148, pushing the existing
# rows 148-153 down to 149-154.
$ws.Rows("148:148").Insert()

$ws.Range("A148").Value = 10
$ws.Range("B148").Value = "Vega Modelo de Temuco"
$ws.Range("C148").Value = "La Araucanía"
$ws.Range("D148").Value = 44516
$ws.Range("E148").Value = 9
$ws.Range("F148").Value = 100112013
$ws.Range("G148").Value = "Alcachofa"
$ws.Range("H148").Value = "Española"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 280
$ws.Range("K148").Value = 300
$ws.Range("L148").Value = 300
$ws.Range("M148").Value = 300
$ws.Range("N148").Value = "$/unidad"
$ws.Range("O148").Value = "Región del Maule"
$ws.Range("P148").Value = 300
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = "Hortaliza"
